$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 65001
$ws.Range("J95").Value = 65001
$ws.Range("L95").Value = 65001
$ws.Range("N95").Value = -70493
$ws.Range("H114").Value = 75000
$ws.Range("J114").Value = 75000
$ws.Range("L114").Value = 75000
$ws.Range("N114").Value = -83678
$ws.Range("H116").Value = 6641.5264
$ws.Range("I116").Value = 6033.5835
$ws.Range("J116").Value = 7683.7144
$ws.Range("K116").Value = 6033.5835
$ws.Range("L116").Value = 7683.7144
$ws.Range("M116").Value = -2591.5835
$ws.Range("N116").Value = -14567.7144
$ws.Range("H137").Value = 3140.8809
$ws.Range("I137").Value = 2316.52
$ws.Range("J137").Value = 4353.1763
$ws.Range("K137").Value = 6949.559999999999
$ws.Range("L137").Value = 13059.5289
$ws.Range("M137").Value = -4399.559999999999
$ws.Range("N137").Value = -18159.5289

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5413386.5
$ws.Range("I2").Value = 8365606.5
$ws.Range("J2").Value = 982.6667
$ws.Range("K2").Value = 8365606.5
$ws.Range("L2").Value = 982.6667
$ws.Range("M2").Value = -8365493.5
$ws.Range("N2").Value = -1208.6667
$ws.Range("H32").Value = 31081.904
$ws.Range("I32").Value = 31327.805
$ws.Range("K32").Value = 31327.805
$ws.Range("M32").Value = -31040.805
$ws.Range("H61").Value = 14499149
$ws.Range("I61").Value = 20837304
$ws.Range("K61").Value = 20837304
$ws.Range("M61").Value = -20837092
$ws.Range("H116").Value = 5413386.5
$ws.Range("I116").Value = 8365606.5
$ws.Range("J116").Value = 982.6667
$ws.Range("K116").Value = 8365606.5
$ws.Range("L116").Value = 982.6667
$ws.Range("M116").Value = -8363312.5
$ws.Range("N116").Value = -5570.6667
$ws.Range("H132").Value = 7148092
$ws.Range("I132").Value = 9527662
$ws.Range("K132").Value = 28582986
$ws.Range("M132").Value = -28580456
$ws.Range("H136").Value = 14499149
$ws.Range("I136").Value = 20837304
$ws.Range("K136").Value = 62511912
$ws.Range("M136").Value = -62509362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5413386.5
$ws.Range("I3").Value = 8365606.5
$ws.Range("J3").Value = 982.6667
$ws.Range("K3").Value = 8365606.5
$ws.Range("L3").Value = 982.6667
$ws.Range("M3").Value = -8365492.5
$ws.Range("N3").Value = -1210.6667
$ws.Range("H86").Value = 1249.3
$ws.Range("I86").Value = 1099.1333
$ws.Range("J86").Value = 1699.8
$ws.Range("K86").Value = 1099.1333
$ws.Range("L86").Value = 1699.8
$ws.Range("M86").Value = 23.86670000000004
$ws.Range("N86").Value = -3945.8
$ws.Range("H89").Value = 1249.3
$ws.Range("I89").Value = 1099.1333
$ws.Range("J89").Value = 1699.8
$ws.Range("K89").Value = 5495.666499999999
$ws.Range("L89").Value = 8499
$ws.Range("M89").Value = 120.3335000000006
$ws.Range("N89").Value = -19731

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 20000
$ws.Range("J28").Value = 20000
$ws.Range("L28").Value = 20000
$ws.Range("N28").Value = -20490
$ws.Range("H31").Value = 4414.5713
$ws.Range("I31").Value = 2632.2307
$ws.Range("K31").Value = 2632.2307
$ws.Range("M31").Value = -2337.2307
$ws.Range("H34").Value = 4414.5713
$ws.Range("I34").Value = 2632.2307
$ws.Range("K34").Value = 2632.2307
$ws.Range("M34").Value = -2430.2307
$ws.Range("H58").Value = 4281.643
$ws.Range("I58").Value = 2685.6428
$ws.Range("K58").Value = 2685.6428
$ws.Range("M58").Value = -2482.6428
$ws.Range("H59").Value = 49362.375
$ws.Range("J59").Value = 52483.168
$ws.Range("L59").Value = 52483.168
$ws.Range("N59").Value = -54773.168
$ws.Range("H62").Value = 9209.066000000001
$ws.Range("I62").Value = 7764.7
$ws.Range("J62").Value = 12097.8
$ws.Range("K62").Value = 7764.7
$ws.Range("L62").Value = 12097.8
$ws.Range("M62").Value = -7140.7
$ws.Range("N62").Value = -13345.8
$ws.Range("H65").Value = 9209.066000000001
$ws.Range("I65").Value = 7764.7
$ws.Range("J65").Value = 12097.8
$ws.Range("K65").Value = 38823.5
$ws.Range("L65").Value = 60489
$ws.Range("M65").Value = -35703.5
$ws.Range("N65").Value = -66729
$ws.Range("H74").Value = 52898.6
$ws.Range("J74").Value = 52898.6
$ws.Range("L74").Value = 52898.6
$ws.Range("N74").Value = -54646.6
$ws.Range("H77").Value = 52898.6
$ws.Range("J77").Value = 52898.6
$ws.Range("L77").Value = 158695.8
$ws.Range("N77").Value = -167431.8
$ws.Range("H132").Value = 23160.355
$ws.Range("I132").Value = 3929.9092
$ws.Range("K132").Value = 11789.7276
$ws.Range("M132").Value = -9259.7276
$ws.Range("H136").Value = 4281.643
$ws.Range("I136").Value = 2685.6428
$ws.Range("K136").Value = 8056.928400000001
$ws.Range("M136").Value = -5506.928400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 22671.666
$ws.Range("J33").Value = 22671.666
$ws.Range("L33").Value = 22671.666
$ws.Range("N33").Value = -23175.666
$ws.Range("H126").Value = 2531.8928
$ws.Range("I126").Value = 2121.1765
$ws.Range("K126").Value = 6363.529500000001
$ws.Range("M126").Value = -3893.529500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 169000.5
$ws.Range("I7").Value = 202300.6
$ws.Range("K7").Value = 202300.6
$ws.Range("M7").Value = -202188.6
$ws.Range("H22").Value = 3242.8538
$ws.Range("I22").Value = 2008.45
$ws.Range("K22").Value = 2008.45
$ws.Range("M22").Value = -1713.45
$ws.Range("H27").Value = 3242.8538
$ws.Range("I27").Value = 2008.45
$ws.Range("K27").Value = 2008.45
$ws.Range("M27").Value = -1901.45
$ws.Range("H68").Value = 4149.875
$ws.Range("I68").Value = 2312.25
$ws.Range("J68").Value = 5987.5
$ws.Range("K68").Value = 2312.25
$ws.Range("L68").Value = 5987.5
$ws.Range("M68").Value = -1563.25
$ws.Range("N68").Value = -7485.5
$ws.Range("H71").Value = 4149.875
$ws.Range("I71").Value = 2312.25
$ws.Range("J71").Value = 5987.5
$ws.Range("K71").Value = 11561.25
$ws.Range("L71").Value = 29937.5
$ws.Range("M71").Value = -7817.25
$ws.Range("N71").Value = -37425.5
$ws.Range("H100").Value = 17860666
$ws.Range("I100").Value = 83335830
$ws.Range("J100").Value = 3802.4546
$ws.Range("K100").Value = 83335830
$ws.Range("L100").Value = 3802.4546
$ws.Range("M100").Value = -83335289
$ws.Range("N100").Value = -4884.4546
$ws.Range("H122").Value = 3009.4443
$ws.Range("I122").Value = 3009.4443
$ws.Range("K122").Value = 9028.332900000001
$ws.Range("M122").Value = -6578.332900000001
$ws.Range("H124").Value = 65065.6
$ws.Range("J124").Value = 65065.6
$ws.Range("L124").Value = 65065.6
$ws.Range("N124").Value = -74885.60000000001
$ws.Range("H126").Value = 169000.5
$ws.Range("I126").Value = 202300.6
$ws.Range("K126").Value = 606901.8
$ws.Range("M126").Value = -604431.8
$ws.Range("H132").Value = 4904.024
$ws.Range("I132").Value = 3831.138
$ws.Range("K132").Value = 11493.414
$ws.Range("M132").Value = -8963.414000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16333.167
$ws.Range("I81").Value = 7799.6
$ws.Range("K81").Value = 15599.2
$ws.Range("M81").Value = -14538.2
$ws.Range("H84").Value = 16333.167
$ws.Range("I84").Value = 7799.6
$ws.Range("K84").Value = 77996
$ws.Range("M84").Value = -72692
$ws.Range("H100").Value = 2681.4
$ws.Range("I100").Value = 802
$ws.Range("K100").Value = 1604
$ws.Range("M100").Value = -1063
$ws.Range("H101").Value = 35198.4
$ws.Range("J101").Value = 35198.4
$ws.Range("L101").Value = 35198.4
$ws.Range("N101").Value = -41688.4
$ws.Range("H132").Value = 4523.8213
$ws.Range("I132").Value = 3959.6667
$ws.Range("J132").Value = 6831.727
$ws.Range("K132").Value = 11879.0001
$ws.Range("L132").Value = 20495.181
$ws.Range("M132").Value = -9349.000100000001
$ws.Range("N132").Value = -25555.181
$ws.Range("H136").Value = 3863503.2
$ws.Range("I136").Value = 6804197
$ws.Range("J136").Value = 3842.1875
$ws.Range("K136").Value = 20412591
$ws.Range("M136").Value = -20410041
$ws.Range("N136").Value = -16626.5625
